$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A47").Value = "NSC"
$ws.Range("A49").Value = "LNCCI"
